# "first 100 abstracts done"
# Mark column A ("status") for the first 100 references (rows 2-101) as
# reviewed. Most entries are plain journal-article citations and get
# "done"; a handful of entries that are databases/web resources rather
# than citable works are instead flagged "X".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A5").Value = "done"
$ws.Range("A6").Value = "X"
$ws.Range("A7:A13").Value = "done"
$ws.Range("A14").Value = "X"
$ws.Range("A15:A20").Value = "done"
$ws.Range("A21").Value = "X"
$ws.Range("A22:A40").Value = "done"
$ws.Range("A41").Value = "X"
$ws.Range("A42:A43").Value = "done"
$ws.Range("A44").Value = "X"
$ws.Range("A45:A77").Value = "done"
$ws.Range("A78").Value = "X"
$ws.Range("A79").Value = "done"
$ws.Range("A80").Value = "X"
$ws.Range("A81:A101").Value = "done"

# Leave the cursor where the author stopped working (row ~102) so the
# saved view reflects where review left off.
$ws.Range("A102").Select()
